$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Waspada"
$ws.Range("B5").Value = "Waspada"
$ws.Range("B6").Value = "Waspada"
$ws.Range("B7").Value = "Waspada"
$ws.Range("B8").Value = "Aman"
$ws.Range("B9").Value = "Waspada"
$ws.Range("B14").Value = "Waspada"
$ws.Range("B15").Value = "Waspada"
$ws.Range("B17").Value = "Waspada"
$ws.Range("B19").Value = "Waspada"
